$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2759
$ws.Range("K3").Value = 2666
$ws.Range("I4").Value = 1789
$ws.Range("K4").Value = 557
$ws.Range("K5").Value = 177
$ws.Range("J6").Value = 11057
$ws.Range("K6").Value = 3311
$ws.Range("I7").Value = 26242
$ws.Range("J7").Value = 29285
$ws.Range("K7").Value = 9470

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 183
$ws.Range("K3").Value = 187
$ws.Range("K7").Value = 625

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 79
$ws.Range("K7").Value = 206

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 102
$ws.Range("K3").Value = 136
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 372

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 105
$ws.Range("J6").Value = 261
$ws.Range("K6").Value = 99
$ws.Range("J7").Value = 901
$ws.Range("K7").Value = 311

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 171

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 71
$ws.Range("K7").Value = 279
$ws.Range("K8").Value = 625
$ws.Range("K9").Value = 35
$ws.Range("K10").Value = 52
$ws.Range("K11").Value = 201
$ws.Range("K15").Value = 96
$ws.Range("K18").Value = 64
$ws.Range("K19").Value = 278
$ws.Range("K20").Value = 217
$ws.Range("K21").Value = 26
$ws.Range("K22").Value = 30
$ws.Range("K23").Value = 82
$ws.Range("K26").Value = 15
$ws.Range("K29").Value = 492
$ws.Range("K31").Value = 109
$ws.Range("K33").Value = 372
$ws.Range("K35").Value = 15
$ws.Range("K36").Value = 109
$ws.Range("J37").Value = 901
$ws.Range("K37").Value = 311
$ws.Range("K40").Value = 21
$ws.Range("K42").Value = 333
$ws.Range("K43").Value = 81
$ws.Range("K46").Value = 20
$ws.Range("K48").Value = 114
$ws.Range("K51").Value = 106
$ws.Range("K52").Value = 262
$ws.Range("K53").Value = 138
$ws.Range("K54").Value = 175
$ws.Range("I63").Value = 202
$ws.Range("K63").Value = 31
$ws.Range("K67").Value = 367
$ws.Range("K72").Value = 43
$ws.Range("K73").Value = 93
$ws.Range("K75").Value = 35
$ws.Range("K76").Value = 140
$ws.Range("K77").Value = 67
$ws.Range("K79").Value = 242
$ws.Range("K83").Value = 206
$ws.Range("K84").Value = 67
$ws.Range("K85").Value = 456
$ws.Range("K86").Value = 59
$ws.Range("K88").Value = 107
$ws.Range("K89").Value = 124
$ws.Range("K90").Value = 87
$ws.Range("K93").Value = 40
$ws.Range("K94").Value = 113
$ws.Range("K95").Value = 154
$ws.Range("K96").Value = 131
$ws.Range("K97").Value = 81
$ws.Range("K98").Value = 53
$ws.Range("K99").Value = 171
$ws.Range("I101").Value = 26242
$ws.Range("J101").Value = 29285
$ws.Range("K101").Value = 9470

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 118
$ws.Range("K4").Value = 23
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 367

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 175

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 131
$ws.Range("K3").Value = 164
$ws.Range("K7").Value = 492

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 73
$ws.Range("K6").Value = 93
$ws.Range("K7").Value = 278

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 106
$ws.Range("K6").Value = 128
$ws.Range("K7").Value = 333

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 79
$ws.Range("K3").Value = 88
$ws.Range("K5").Value = 8
$ws.Range("K7").Value = 242

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 93
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K4").Value = 9
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 21
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K4").Value = 20
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K4").Value = 8
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 167
$ws.Range("K3").Value = 156
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 456

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 262
